$d = $word.ActiveDocument

$replacements = @(
    @("832÷2=416, 0", "739÷6=123, 1"),
    @("890÷2=445, 0", "504÷7=72, 0"),
    @("845÷7=120, 5", "629÷9=69, 8"),
    @("870÷6=145, 0", "309÷6=51, 3"),
    @("807÷8=100, 7", "934÷4=233, 2"),
    @("113÷7=16, 1",  "486÷7=69, 3"),
    @("366÷4=91, 2",  "468÷5=93, 3"),
    @("556÷3=185, 1", "988÷7=141, 1"),
    @("708÷6=118, 0", "967÷7=138, 1"),
    @("319÷5=63, 4",  "764÷8=95, 4"),
    @("164÷2=82, 0",  "747÷5=149, 2"),
    @("504÷9=56, 0",  "908÷5=181, 3"),
    @("942÷3=314, 0", "277÷7=39, 4"),
    @("349÷5=69, 4",  "391÷5=78, 1"),
    @("131÷4=32, 3",  "361÷6=60, 1"),
    @("611÷9=67, 8",  "638÷6=106, 2"),
    @("437÷4=109, 1", "795÷5=159, 0"),
    @("982÷9=109, 1", "433÷6=72, 1"),
    @("320÷2=160, 0", "420÷2=210, 0"),
    @("400÷7=57, 1",  "664÷8=83, 0"),
    @("779÷7=111, 2", "777÷5=155, 2"),
    @("396÷5=79, 1",  "301÷2=150, 1"),
    @("991÷2=495, 1", "282÷7=40, 2"),
    @("938÷6=156, 2", "913÷6=152, 1"),
    @("705÷8=88, 1",  "303÷5=60, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
